# Weekly LDLC price-history snapshot: insert a new timestamp column right
# before the "nom" / "url_produit" columns, pushing them one column to the
# right (CM -> CN, CN -> CO). The new column carries the latest known
# value forward from the previous snapshot column (same pattern already
# used by every prior snapshot column in this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "nom" currently lives in column CM (93 -> but as letters: CM).
# Insert a fresh blank column there; Excel shifts CM->CN and CN->CO.
$ws.Columns("CM").Insert()

# Header for the freshly inserted snapshot column.
$ws.Range("CM1").Value = "2026-01-31 19:13:20"

# CL is the previous (most recent) snapshot column (column index 90).
# CM is the new snapshot column (column index 91), just inserted blank.
$lastCol = 90
$newCol = 91
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $prevValue = $ws.Cells.Item($r, $lastCol).Value()
    if ($null -ne $prevValue -and $prevValue -ne "") {
        $ws.Cells.Item($r, $newCol).Value = $prevValue
    }
}
